# "avoir infos du bouton 1"
# Turn the MDN URL in the Webographie section into a real clickable
# hyperlink, and split the paragraph so the trailing (hidden) _GoBack
# bookmark ends up alone in its own paragraph - exactly what Word does
# when you place the cursor right after the URL and press Enter while
# the URL autoformats into a hyperlink.

$d = $word.ActiveDocument

# Locate the plain-text URL run.
$urlRange = $d.Content.Duplicate
$urlRange.Find.Execute(
    "https://developer.mozilla.org/fr/docs/Web/HTML/Element/Input/range",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Turn it into a hyperlink (TextToDisplay keeps the visible text the same
# as the target address, exactly like the URL-autoformat feature does).
$link = $d.Hyperlinks.Add($urlRange, $urlRange.Text, "", "", $urlRange.Text)
$link.Range.Style = "Lienhypertexte"

# Split the paragraph right after the hyperlink so that the trailing
# _GoBack bookmark (and the paragraph mark that used to follow the URL)
# moves into its own, new paragraph that keeps the same style.
$splitPoint = $d.Range($link.Range.End, $link.Range.End)
$splitPoint.InsertBefore([char]13)
